$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values for columns B:E
$ws.Range("B2").Value = 234.40574972987804
$ws.Range("C2").Value = 172.98601574991676
$ws.Range("D2").Value = 233.87009906794606
$ws.Range("E2").Value = 173.78653927819192

# Update row 3 values for columns B:E
$ws.Range("B3").Value = 214.8102170230153
$ws.Range("C3").Value = 171.66530788094812
$ws.Range("D3").Value = 209.44538012402467
$ws.Range("E3").Value = 176.8876067736812

# Update the selected range on the sheet to reflect the new active selection
$ws.Range("B1:E3").Select()
